$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the exercise link and slides link for session 04 (row 5)
# Order matters for shared-strings table index assignment: exercise string must
# be registered before the slides string so the new entries land at indices 31/32.
$ws.Range("F5").Value = "exercises/e04.html"
$ws.Range("E5").Value = "slides/slides.html#/session-04-structuring-a-heterogeneous-field-the-basics-of-markdown-and-github"

# Update the active cell selection to E5 as reflected in the saved file
$ws.Range("E5").Select()
